# Applies the "add showifdef from teaching tuesday" edit:
#  1. Removes the spell-check proofErr wrapper around the "Defendant"
#     run in the Defendant/defendant_experts table cell (2nd "Defendant:"
#     cell - the 1st, in the witnesses table, already has no proofErr).
#  2. Wraps {{ users[0].signature }} in showifdef(...) ->
#     {{ showifdef('users[0].signature')}}
#  3. Wraps {{ signature_date }} in showifdef(...) ->
#     {{ showifdef('signature_date') }}

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument
$paras = $d.Paragraphs

# ---------------------------------------------------------------------
# Helper: Paragraph.Range.Text carries trailing control marks - a
# paragraph mark (13) and, for the last paragraph in a table cell, a
# cell mark (7) as well. Strip those off before comparing text.
# ---------------------------------------------------------------------
function Normalize-ParaText($t) {
    return $t.TrimEnd([char]13, [char]7)
}

# ---------------------------------------------------------------------
# Helper: find a paragraph whose trimmed text equals $text, optionally
# constrained by the trimmed text of the following paragraph.
# ---------------------------------------------------------------------
function Find-ParagraphByText($paras, $text, $nextText) {
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = Normalize-ParaText $p.Range.Text
        if ($t -eq $text) {
            if ($nextText -ne $null) {
                if (($i + 1) -le $paras.Count) {
                    $nt = Normalize-ParaText $paras.Item($i + 1).Range.Text
                    if ($nt -ne $nextText) {
                        continue
                    }
                } else {
                    continue
                }
            }
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. Defendant: cell (defendant_experts table) - drop proofErr wrapper
# ---------------------------------------------------------------------
$pDefendant = Find-ParagraphByText $paras "Defendant:" "{{ defendant_experts }}"
if ($pDefendant -eq $null) {
    throw "Could not locate the 'Defendant:' paragraph (defendant_experts cell)."
}
$xmlDefendant = '<w:p ' + $wordNs + '><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Defendant</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
$pDefendant.Range.InsertXML($xmlDefendant)

# ---------------------------------------------------------------------
# 2. {{ users[0].signature }} -> {{ showifdef('users[0].signature')}}
# ---------------------------------------------------------------------
$pSig = Find-ParagraphByText $paras "{{ users[0].signature }}`t" $null
if ($pSig -eq $null) {
    throw "Could not locate the '{{ users[0].signature }}' paragraph."
}
$xmlSig = '<w:p ' + $wordNs + '>' +
  '<w:pPr><w:keepNext/><w:keepLines/><w:tabs><w:tab w:val="left" w:pos="9360"/></w:tabs><w:ind w:left="5040"/><w:contextualSpacing/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>showifdef</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(''</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>users[0].signature</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>'')</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>}}</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:tab/></w:r>' +
  '</w:p>'
$pSig.Range.InsertXML($xmlSig)

# ---------------------------------------------------------------------
# 3. Dated: {{ signature_date }} -> Dated: {{ showifdef('signature_date') }}
# ---------------------------------------------------------------------
$pDate = Find-ParagraphByText $paras "Dated: {{ signature_date }}" $null
if ($pDate -eq $null) {
    throw "Could not locate the 'Dated: {{ signature_date }}' paragraph."
}
$xmlDate = '<w:p ' + $wordNs + '>' +
  '<w:pPr><w:keepLines/><w:tabs><w:tab w:val="left" w:pos="3560"/><w:tab w:val="left" w:pos="4320"/></w:tabs></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Dated: </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>showifdef</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(''</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>signature_date</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">'') </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>}}</w:t></w:r>' +
  '</w:p>'
$pDate.Range.InsertXML($xmlDate)

Write-Output "Edits applied."
